$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-06-10 Monday"; new = "2024-06-11 Tuesday"},
    @{old = "420÷2="; new = "196÷9="},
    @{old = "687÷8="; new = "947÷7="},
    @{old = "240÷6="; new = "226÷2="},
    @{old = "623÷6="; new = "370÷7="},
    @{old = "391÷6="; new = "916÷4="},
    @{old = "253÷2="; new = "914÷6="},
    @{old = "719÷9="; new = "229÷8="},
    @{old = "317÷3="; new = "681÷8="},
    @{old = "530÷8="; new = "382÷5="},
    @{old = "859÷9="; new = "231÷4="},
    @{old = "803÷2="; new = "576÷6="},
    @{old = "969÷8="; new = "334÷7="},
    @{old = "157÷3="; new = "559÷3="},
    @{old = "658÷6="; new = "271÷7="},
    @{old = "880÷8="; new = "249÷9="},
    @{old = "271÷8="; new = "846÷5="},
    @{old = "265÷8="; new = "570÷7="},
    @{old = "682÷9="; new = "972÷3="},
    @{old = "155÷3="; new = "341÷8="},
    @{old = "672÷5="; new = "172÷4="},
    @{old = "317÷6="; new = "606÷9="},
    @{old = "544÷8="; new = "634÷7="},
    @{old = "352÷8="; new = "993÷3="},
    @{old = "215÷2="; new = "936÷4="},
    @{old = "408÷8="; new = "863÷5="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
